# Apply review updates to the "Review" sheet:
#  - Row 4 (the "Document status section..." review point) date corrected
#    from 24/1/2021 to 24/1/2020 (aligning with the other rows for the
#    same document / review date group).
#  - Several review points' status flipped from "Open" to "Closed" as the
#    CYRS/HSI structural issues they tracked have been resolved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the date on review point in row 4 (was mistakenly 24/1/2021).
$ws.Range("A4").Value = "24/1/2020"

# Mark these review points as Closed now that the CYRS/HSI structure was reworked.
$ws.Range("E3").Value = "Closed"
$ws.Range("E4").Value = "Closed"
$ws.Range("E6").Value = "Closed"
$ws.Range("E11").Value = "Closed"

# Leave the selection on the last updated cell, matching where review left off.
$ws.Range("E11").Select()
